$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("210:210").Insert()
$ws.Range("A210").Value = 5
$ws.Range("B210").Value = "Macroferia Regional de Talca"
$ws.Range("C210").Value = "Maule"
$ws.Range("D210").Value = 44627
$ws.Range("D210").NumberFormat = $ws.Range("D211").NumberFormat
$ws.Range("E210").Value = 7
$ws.Range("F210").Value = 100112023
$ws.Range("G210").Value = "Brócoli"
$ws.Range("H210").Value = "Sin especificar"
$ws.Range("I210").Value = "Primera"
$ws.Range("J210").Value = 4000
$ws.Range("K210").Value = 600
$ws.Range("L210").Value = 600
$ws.Range("M210").Value = 600
$ws.Range("N210").Value = "$/unidad"
$ws.Range("O210").Value = "Región del Maule"
$ws.Range("P210").Value = 600
$ws.Range("Q210").Value = 1
$ws.Range("R210").Value = "Hortaliza"
